$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 97
$ws.Range("I12").Value = 85
$ws.Range("J12").Value = 115
$ws.Range("K12").Value = 85
$ws.Range("L12").Value = 115
$ws.Range("M12").Value = 85
$ws.Range("N12").Value = -455
$ws.Range("H98").Value = 1447.8125
$ws.Range("J98").Value = 833
$ws.Range("L98").Value = 833
$ws.Range("N98").Value = -3829
$ws.Range("H116").Value = 4508.25
$ws.Range("I116").Value = 5187.222
$ws.Range("J116").Value = 3952.7273
$ws.Range("K116").Value = 5187.222
$ws.Range("L116").Value = 3952.7273
$ws.Range("M116").Value = -1745.222
$ws.Range("N116").Value = -10836.7273
$ws.Range("H122").Value = 1447.8125
$ws.Range("J122").Value = 833
$ws.Range("L122").Value = 2499
$ws.Range("N122").Value = -7399
$ws.Range("H131").Value = 1349.9286
$ws.Range("I131").Value = 598.5454999999999
$ws.Range("J131").Value = 4105
$ws.Range("K131").Value = 1795.6365
$ws.Range("L131").Value = 12315
$ws.Range("M131").Value = 3244.3635
$ws.Range("N131").Value = -22395
$ws.Range("H132").Value = 1416.6327
$ws.Range("I132").Value = 1448.881
$ws.Range("J132").Value = 1223.1428
$ws.Range("K132").Value = 4346.643
$ws.Range("L132").Value = 3669.4284
$ws.Range("M132").Value = -1816.643
$ws.Range("N132").Value = -8729.428400000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1789.9375
$ws.Range("I61").Value = 1552.0834
$ws.Range("J61").Value = 2503.5
$ws.Range("K61").Value = 1552.0834
$ws.Range("L61").Value = 2503.5
$ws.Range("M61").Value = -1340.0834
$ws.Range("N61").Value = -2927.5
$ws.Range("H109").Value = 29800
$ws.Range("J109").Value = 29800
$ws.Range("L109").Value = 29800
$ws.Range("N109").Value = -32574
$ws.Range("H122").Value = 1599.6666
$ws.Range("I122").Value = 919.6
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 2758.8
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -308.8000000000002
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 2000.3715
$ws.Range("I132").Value = 1338.05
$ws.Range("K132").Value = 4014.15
$ws.Range("M132").Value = -1484.15
$ws.Range("H136").Value = 1789.9375
$ws.Range("I136").Value = 1552.0834
$ws.Range("J136").Value = 2503.5
$ws.Range("K136").Value = 4656.2502
$ws.Range("L136").Value = 7510.5
$ws.Range("M136").Value = -2106.2502
$ws.Range("N136").Value = -12610.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 35000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 35000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 35000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -36996
$ws.Range("H84").Value = 35000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 35000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 105000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -114984
$ws.Range("H132").Value = 1797.4615
$ws.Range("I132").Value = 1232.8948
$ws.Range("J132").Value = 3329.8572
$ws.Range("K132").Value = 3698.6844
$ws.Range("L132").Value = 9989.571599999999
$ws.Range("M132").Value = -1168.6844
$ws.Range("N132").Value = -15049.5716
$ws.Range("H134").Value = 55556810
$ws.Range("I134").Value = 1110.7142
$ws.Range("J134").Value = 250001740
$ws.Range("K134").Value = 3332.1426
$ws.Range("L134").Value = 750005220
$ws.Range("M134").Value = -797.1425999999997
$ws.Range("N134").Value = -750010290

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 683.3333
$ws.Range("J92").Value = 725
$ws.Range("L92").Value = 2175
$ws.Range("N92").Value = -4671
$ws.Range("H136").Value = 61739.35
$ws.Range("I136").Value = 126177.375
$ws.Range("J136").Value = 4461.1113
$ws.Range("K136").Value = 378532.125
$ws.Range("L136").Value = 13383.3339
$ws.Range("M136").Value = -373432.125
$ws.Range("N136").Value = -23583.3339

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2307.4
$ws.Range("I102").Value = 2196.75
$ws.Range("K102").Value = 2196.75
$ws.Range("M102").Value = -574.75
$ws.Range("H126").Value = 5043
$ws.Range("I126").Value = 4050.1667
$ws.Range("K126").Value = 12150.5001
$ws.Range("M126").Value = -9680.500100000001
$ws.Range("H132").Value = 2384.0715
$ws.Range("I132").Value = 1135
$ws.Range("J132").Value = 4049.5
$ws.Range("K132").Value = 3405
$ws.Range("L132").Value = 12148.5
$ws.Range("M132").Value = -875
$ws.Range("N132").Value = -17208.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3100
$ws.Range("I7").Value = 3420
$ws.Range("J7").Value = 2300
$ws.Range("K7").Value = 3420
$ws.Range("L7").Value = 2300
$ws.Range("M7").Value = -3308
$ws.Range("N7").Value = -2524
$ws.Range("H22").Value = 737.75
$ws.Range("H27").Value = 737.75
$ws.Range("H46").Value = 1480.2
$ws.Range("I46").Value = 1657.2858
$ws.Range("K46").Value = 1657.2858
$ws.Range("M46").Value = -1469.2858
$ws.Range("H126").Value = 3100
$ws.Range("I126").Value = 3420
$ws.Range("J126").Value = 2300
$ws.Range("K126").Value = 10260
$ws.Range("L126").Value = 6900
$ws.Range("M126").Value = -7790
$ws.Range("N126").Value = -11840
$ws.Range("H132").Value = 1890.1628
$ws.Range("I132").Value = 1211.5714
$ws.Range("J132").Value = 3156.8667
$ws.Range("K132").Value = 3634.7142
$ws.Range("L132").Value = 9470.6001
$ws.Range("M132").Value = -1104.7142
$ws.Range("N132").Value = -14530.6001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 970.8788
$ws.Range("I107").Value = 1017.7727
$ws.Range("J107").Value = 877.0909
$ws.Range("K107").Value = 3053.3181
$ws.Range("L107").Value = 2631.2727
$ws.Range("M107").Value = -1133.3181
$ws.Range("N107").Value = -6471.2727
$ws.Range("H132").Value = 1375.15
$ws.Range("I132").Value = 1112.3077
$ws.Range("K132").Value = 3336.9231
$ws.Range("M132").Value = -806.9231
$ws.Range("H136").Value = 356.60526
$ws.Range("I136").Value = 334.76666
$ws.Range("J136").Value = 438.5
$ws.Range("K136").Value = 1004.29998
$ws.Range("L136").Value = 1315.5
$ws.Range("M136").Value = 1545.70002
$ws.Range("N136").Value = -6415.5
